$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for one additional data row above the "Total" row.
# Row 15 is currently empty (unused) and the "Total" row is row 16;
# inserting a row at 16 shifts "Total" (and the hours row below it) down by one,
# giving us two free rows (15 and 16) for the new entries.
$ws.Rows("16:16").Insert()

# New entries (set text cells first, in the same order they were
# authored, so new shared-string entries come out in the expected order)
$ws.Range("C16").Value = "Code: Extend evaluation metrics, adjust lr, adjust reward weighting, troubleshooting"
$ws.Range("C15").Value = "Code: adjust reward"
$ws.Range("A15").Value = "31.08."
$ws.Range("A16").Value = "02.09."
$ws.Range("B15").Value = 75
$ws.Range("B16").Value = 520

# Update the Total formulas (now at rows 17/18) to include the new rows up to B16
$ws.Range("B17").Formula = '=_xlfn.CONCAT(SUM(B4:B16)," min")'
$ws.Range("B18").Formula = '=_xlfn.CONCAT("~",ROUND(SUM(B4:B16) / 60, 2)," h")'

# Match the selection shown in the saved workbook
$ws.Range("C16").Select()
